$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Backlog row 3 (the "login / connect" story) is now done: move its
# Statut from "A faire" to "Fais".
$ws.Range("E3").Value = "Fais"

# Reflect the last place the user was working before saving.
$ws.Range("B8").Select()
